$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MPAA")

$ws.Range("D43").Value = 157500
$ws.Range("D44").Value = 329900
$ws.Range("D45").Value = 20100
$ws.Range("D46").Value = 287600
$ws.Range("D47").Value = 228300
$ws.Range("D49").Value = 10100
$ws.Range("D52").Value = 309200
$ws.Range("D54").Value = 552400
$ws.Range("D59").Value = 83500
$ws.Range("D60").Value = 197300
$ws.Range("D66").Value = 265500
$ws.Range("D72").Value = 78500
$ws.Range("D76").Value = 286900
